$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G4 with new shared string text
$ws.Range("G4").Value = "On-mouse events, reorganizing code, using SDL_mixer for sound."

# Update F5 and G5 with hour values
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 2

# Update B6 with new text (this also causes row 6 to auto-fit taller)
$ws.Range("B6").Value = "Making buttons play sounds. Using Timing as a counter and to count frames. Capping fps. Collision detections (Squared, Circular and Per-Pixel)."
$ws.Rows.Item(6).RowHeight = 100.8

# Update selection to B7
$ws.Range("B7").Select()
